$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.282.04"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.493.58"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.483"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").Value = "  +6.36%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "4.089.34"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "3.491.90"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "64.225.36"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.577"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").Value = "3.634.35"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("D34").Value = "3.524.17"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.99%  "
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "2.381.88"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("E51").Value = "  -2.08%  "
